$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3657.625
$ws.Range("J40").Value = 3840.3333
$ws.Range("L40").Value = 3840.3333
$ws.Range("N40").Value = -4190.3333
$ws.Range("H43").Value = 3667
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3138
$ws.Range("H55").Value = 171.26315
$ws.Range("I55").Value = 67.2
$ws.Range("J55").Value = 286.8889
$ws.Range("K55").Value = 67.2
$ws.Range("L55").Value = 286.8889
$ws.Range("M55").Value = 146.8
$ws.Range("N55").Value = -714.8888999999999
$ws.Range("H86").Value = 5500
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377
$ws.Range("H89").Value = 5500
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884
$ws.Range("H97").Value = 866
$ws.Range("J97").Value = 866
$ws.Range("L97").Value = 2598
$ws.Range("N97").Value = -3590
$ws.Range("H115").Value = 627.5714
$ws.Range("I115").Value = 627.5714
$ws.Range("K115").Value = 1882.7142
$ws.Range("M115").Value = -315.7142000000001
$ws.Range("H135").Value = 2758.1667
$ws.Range("I135").Value = 2631.9375
$ws.Range("J135").Value = 3768
$ws.Range("K135").Value = 23687.4375
$ws.Range("L135").Value = 33912
$ws.Range("M135").Value = -21152.4375
$ws.Range("N135").Value = -38982

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3650.8518
$ws.Range("I32").Value = 3732.2292
$ws.Range("K32").Value = 3732.2292
$ws.Range("M32").Value = -3445.2292
$ws.Range("H47").Value = 39999.4
$ws.Range("J47").Value = 39999.4
$ws.Range("L47").Value = 39999.4
$ws.Range("N47").Value = -41449.4
$ws.Range("H61").Value = 2574.3
$ws.Range("I61").Value = 2440.647
$ws.Range("J61").Value = 3331.6667
$ws.Range("K61").Value = 2440.647
$ws.Range("L61").Value = 3331.6667
$ws.Range("M61").Value = -2228.647
$ws.Range("N61").Value = -3755.6667
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240
$ws.Range("H74").Value = 1406
$ws.Range("I74").Value = 1491.0286
$ws.Range("J74").Value = 980.8570999999999
$ws.Range("K74").Value = 1491.0286
$ws.Range("L74").Value = 980.8570999999999
$ws.Range("M74").Value = -617.0286000000001
$ws.Range("N74").Value = -2728.8571
$ws.Range("H77").Value = 1406
$ws.Range("I77").Value = 1491.0286
$ws.Range("J77").Value = 980.8570999999999
$ws.Range("K77").Value = 7455.143
$ws.Range("L77").Value = 4904.2855
$ws.Range("M77").Value = -3087.143
$ws.Range("N77").Value = -13640.2855
$ws.Range("H102").Value = 1908.4
$ws.Range("I102").Value = 1908.4
$ws.Range("K102").Value = 1908.4
$ws.Range("M102").Value = -286.4000000000001
$ws.Range("H132").Value = 1671.7878
$ws.Range("I132").Value = 1671.7878
$ws.Range("K132").Value = 5015.3634
$ws.Range("M132").Value = -2485.3634
$ws.Range("H136").Value = 2574.3
$ws.Range("I136").Value = 2440.647
$ws.Range("J136").Value = 3331.6667
$ws.Range("K136").Value = 7321.941
$ws.Range("L136").Value = 9995.000100000001
$ws.Range("M136").Value = -4771.941
$ws.Range("N136").Value = -15095.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4172.2144
$ws.Range("I86").Value = 2934.6667
$ws.Range("K86").Value = 2934.6667
$ws.Range("M86").Value = -1811.6667
$ws.Range("H89").Value = 4172.2144
$ws.Range("I89").Value = 2934.6667
$ws.Range("K89").Value = 14673.3335
$ws.Range("M89").Value = -9057.333500000001
$ws.Range("H134").Value = 3164.1304
$ws.Range("I134").Value = 2253.9412
$ws.Range("J134").Value = 5743
$ws.Range("K134").Value = 6761.823600000001
$ws.Range("L134").Value = 17229
$ws.Range("M134").Value = -4226.823600000001
$ws.Range("N134").Value = -22299
$ws.Range("H140").Value = 79090.63
$ws.Range("J140").Value = 79090.63
$ws.Range("L140").Value = 79090.63
$ws.Range("N140").Value = -89450.63

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1399.3903
$ws.Range("I31").Value = 1290.0605
$ws.Range("K31").Value = 1290.0605
$ws.Range("M31").Value = -995.0605
$ws.Range("H34").Value = 1399.3903
$ws.Range("I34").Value = 1290.0605
$ws.Range("K34").Value = 1290.0605
$ws.Range("M34").Value = -1088.0605
$ws.Range("H99").Value = 4069.36
$ws.Range("I99").Value = 3768.7
$ws.Range("K99").Value = 3768.7
$ws.Range("M99").Value = -2270.7
$ws.Range("H103").Value = 20000
$ws.Range("I103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("K103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("M103").Value = -18828
$ws.Range("N103").Value = -22344
$ws.Range("H126").Value = 4069.36
$ws.Range("I126").Value = 3768.7
$ws.Range("K126").Value = 11306.1
$ws.Range("M126").Value = -8836.099999999999
$ws.Range("H132").Value = 2610.5925
$ws.Range("I132").Value = 1750.6666
$ws.Range("K132").Value = 5251.9998
$ws.Range("M132").Value = -2721.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 23.75
$ws.Range("J33").Value = 18.333334
$ws.Range("L33").Value = 110.000004
$ws.Range("N33").Value = -676.000004
$ws.Range("H94").Value = 3026.1428
$ws.Range("I94").Value = 296
$ws.Range("K94").Value = 888
$ws.Range("M94").Value = -212
$ws.Range("H129").Value = 1235.7858
$ws.Range("I129").Value = 600.25
$ws.Range("J129").Value = 2083.1667
$ws.Range("K129").Value = 1800.75
$ws.Range("L129").Value = 6249.500100000001
$ws.Range("M129").Value = 3199.25
$ws.Range("N129").Value = -16249.5001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 24999.5
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H44").Value = 21974
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H52").Value = 37826.332
$ws.Range("I52").Value = 37030
$ws.Range("J52").Value = 37985.6
$ws.Range("K52").Value = 37030
$ws.Range("L52").Value = 37985.6
$ws.Range("M52").Value = -36771
$ws.Range("N52").Value = -38503.6
$ws.Range("H80").Value = 81627.47
$ws.Range("I80").Value = 105061.18
$ws.Range("J80").Value = 17184.75
$ws.Range("K80").Value = 105061.18
$ws.Range("L80").Value = 17184.75
$ws.Range("M80").Value = -104063.18
$ws.Range("N80").Value = -19180.75
$ws.Range("H83").Value = 81627.47
$ws.Range("I83").Value = 105061.18
$ws.Range("J83").Value = 17184.75
$ws.Range("K83").Value = 525305.8999999999
$ws.Range("L83").Value = 85923.75
$ws.Range("M83").Value = -520313.8999999999
$ws.Range("N83").Value = -95907.75
$ws.Range("H102").Value = 3386.158
$ws.Range("I102").Value = 2781.2144
$ws.Range("K102").Value = 2781.2144
$ws.Range("M102").Value = -1159.2144

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 25833.334
$ws.Range("J53").Value = 28000
$ws.Range("L53").Value = 28000
$ws.Range("N53").Value = -29036

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 22500
$ws.Range("I43").Value = 22500
$ws.Range("K43").Value = 22500
$ws.Range("M43").Value = -22351
$ws.Range("H45").Value = 14025.556
$ws.Range("I45").Value = 7967
$ws.Range("J45").Value = 14782.875
$ws.Range("K45").Value = 7967
$ws.Range("L45").Value = 14782.875
$ws.Range("M45").Value = -7476
$ws.Range("N45").Value = -15764.875
$ws.Range("H53").Value = 23831.666
$ws.Range("I53").Value = 1495
$ws.Range("J53").Value = 35000
$ws.Range("K53").Value = 1495
$ws.Range("L53").Value = 35000
$ws.Range("M53").Value = -888
$ws.Range("N53").Value = -36214
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 3968.75
$ws.Range("I81").Value = 2115.5
$ws.Range("J81").Value = 5822
$ws.Range("K81").Value = 4231
$ws.Range("L81").Value = 11644
$ws.Range("M81").Value = -3170
$ws.Range("N81").Value = -13766
$ws.Range("H84").Value = 3968.75
$ws.Range("I84").Value = 2115.5
$ws.Range("J84").Value = 5822
$ws.Range("K84").Value = 21155
$ws.Range("L84").Value = 58220
$ws.Range("M84").Value = -15851
$ws.Range("N84").Value = -68828
$ws.Range("H122").Value = 3217.875
$ws.Range("I122").Value = 3032.1667
$ws.Range("J122").Value = 3775
$ws.Range("K122").Value = 9096.500100000001
$ws.Range("L122").Value = 11325
$ws.Range("M122").Value = -6646.500100000001
$ws.Range("N122").Value = -16225
$ws.Range("H126").Value = 3876.389
$ws.Range("I126").Value = 3986.4707
$ws.Range("J126").Value = 2005
$ws.Range("K126").Value = 11959.4121
$ws.Range("L126").Value = 6015
$ws.Range("M126").Value = -9489.4121
$ws.Range("N126").Value = -10955
$ws.Range("H132").Value = 1927.3529
$ws.Range("I132").Value = 1451
$ws.Range("K132").Value = 4353
$ws.Range("M132").Value = -1823
